$wb = $excel.ActiveWorkbook
$trials = $wb.Worksheets.Item("Trials")
$summary = $wb.Worksheets.Item("Summary")

$trials.Cells.Item(102, 2).Value = 10886
$trials.Cells.Item(102, 3).Value = 12.69341087341309
$trials.Cells.Item(102, 4).Value = 1

$trials.Cells.Item(103, 2).Value = 4741
$trials.Cells.Item(103, 3).Value = 2.147418975830078

$trials.Cells.Item(104, 2).Value = 5075
$trials.Cells.Item(104, 3).Value = 3.13793683052063

$trials.Cells.Item(105, 2).Value = 12309
$trials.Cells.Item(105, 3).Value = 15.35720777511597

$trials.Cells.Item(106, 2).Value = 8695
$trials.Cells.Item(106, 3).Value = 9.823531150817871

$trials.Cells.Item(107, 2).Value = 6545
$trials.Cells.Item(107, 3).Value = 6.852670669555664

$trials.Cells.Item(108, 2).Value = 2781
$trials.Cells.Item(108, 3).Value = 1.830486059188843

$trials.Cells.Item(109, 2).Value = 4635
$trials.Cells.Item(109, 3).Value = 3.624930620193481

$trials.Cells.Item(110, 2).Value = 10072
$trials.Cells.Item(110, 3).Value = 13.06898617744446

$trials.Cells.Item(111, 2).Value = 10213
$trials.Cells.Item(111, 3).Value = 13.95496845245361

$trials.Cells.Item(112, 2).Value = 23397
$trials.Cells.Item(112, 3).Value = 60.00241875648499
$trials.Cells.Item(112, 4).Value = 0

$trials.Cells.Item(113, 2).Value = 6675
$trials.Cells.Item(113, 3).Value = 5.696521759033203

$trials.Cells.Item(114, 2).Value = 11804
$trials.Cells.Item(114, 3).Value = 13.74366140365601

$trials.Cells.Item(115, 2).Value = 3195
$trials.Cells.Item(115, 3).Value = 1.783372640609741

$trials.Cells.Item(116, 2).Value = 3905
$trials.Cells.Item(116, 3).Value = 2.366245031356812

$trials.Cells.Item(117, 2).Value = 8602
$trials.Cells.Item(117, 3).Value = 8.000714540481567

$trials.Cells.Item(118, 2).Value = 14408
$trials.Cells.Item(118, 3).Value = 20.44626379013062
$trials.Cells.Item(118, 4).Value = 1

$trials.Cells.Item(119, 2).Value = 12037
$trials.Cells.Item(119, 3).Value = 14.95835471153259

$trials.Cells.Item(120, 2).Value = 3504
$trials.Cells.Item(120, 3).Value = 1.999760150909424

$trials.Cells.Item(121, 2).Value = 745
$trials.Cells.Item(121, 3).Value = 0.2582452297210693

$trials.Cells.Item(122, 2).Value = 3727
$trials.Cells.Item(122, 3).Value = 1.746082305908203
$trials.Cells.Item(122, 4).Value = 1

$trials.Cells.Item(123, 2).Value = 26106
$trials.Cells.Item(123, 3).Value = 60.00500917434692
$trials.Cells.Item(123, 4).Value = 0

$trials.Cells.Item(124, 2).Value = 10517
$trials.Cells.Item(124, 3).Value = 11.65313529968262

$trials.Cells.Item(125, 2).Value = 25214
$trials.Cells.Item(125, 3).Value = 60.00290274620056

$trials.Cells.Item(126, 2).Value = 4303
$trials.Cells.Item(126, 3).Value = 3.007611513137817

$trials.Cells.Item(127, 2).Value = 1262
$trials.Cells.Item(127, 3).Value = 0.5637445449829102

$trials.Cells.Item(128, 2).Value = 5866
$trials.Cells.Item(128, 3).Value = 4.028760433197021

$trials.Cells.Item(129, 2).Value = 3629
$trials.Cells.Item(129, 3).Value = 2.166686773300171

$trials.Cells.Item(130, 2).Value = 3641
$trials.Cells.Item(130, 3).Value = 2.227350473403931

$trials.Cells.Item(131, 2).Value = 9169
$trials.Cells.Item(131, 3).Value = 10.29142379760742

$trials.Cells.Item(132, 2).Value = 2288
$trials.Cells.Item(132, 3).Value = 0.736778736114502
$trials.Cells.Item(132, 4).Value = 1

$trials.Cells.Item(133, 2).Value = 7838
$trials.Cells.Item(133, 3).Value = 6.269149780273438

$trials.Cells.Item(134, 2).Value = 1147
$trials.Cells.Item(134, 3).Value = 0.3828873634338379

$trials.Cells.Item(135, 2).Value = 2492
$trials.Cells.Item(135, 3).Value = 1.721743106842041

$trials.Cells.Item(136, 2).Value = 14912
$trials.Cells.Item(136, 3).Value = 25.4238498210907

$trials.Cells.Item(137, 2).Value = 1029
$trials.Cells.Item(137, 3).Value = 0.3962762355804443

$trials.Cells.Item(138, 2).Value = 12974
$trials.Cells.Item(138, 3).Value = 17.48500752449036

$trials.Cells.Item(139, 2).Value = 25555
$trials.Cells.Item(139, 3).Value = 60.00132536888123
$trials.Cells.Item(139, 4).Value = 0

$trials.Cells.Item(140, 2).Value = 5826
$trials.Cells.Item(140, 3).Value = 4.775183916091919

$trials.Cells.Item(141, 2).Value = 3509
$trials.Cells.Item(141, 3).Value = 1.973501443862915

$trials.Cells.Item(142, 2).Value = 4051
$trials.Cells.Item(142, 3).Value = 2.434374809265137

$trials.Cells.Item(143, 2).Value = 2236
$trials.Cells.Item(143, 3).Value = 0.987722635269165

$trials.Cells.Item(144, 2).Value = 8963
$trials.Cells.Item(144, 3).Value = 8.611688852310181

$trials.Cells.Item(145, 2).Value = 5651
$trials.Cells.Item(145, 3).Value = 3.819227695465088

$trials.Cells.Item(146, 2).Value = 15898
$trials.Cells.Item(146, 3).Value = 24.02742671966553

$trials.Cells.Item(147, 2).Value = 3440
$trials.Cells.Item(147, 3).Value = 1.893076181411743

$trials.Cells.Item(148, 2).Value = 9353
$trials.Cells.Item(148, 3).Value = 9.305517911911011

$trials.Cells.Item(149, 2).Value = 16540
$trials.Cells.Item(149, 3).Value = 27.72230744361877
$trials.Cells.Item(149, 4).Value = 1

$trials.Cells.Item(150, 2).Value = 3702
$trials.Cells.Item(150, 3).Value = 2.255652904510498

$trials.Cells.Item(151, 2).Value = 3764
$trials.Cells.Item(151, 3).Value = 2.236135721206665
$trials.Cells.Item(151, 4).Value = 1

$trials.Cells.Item(152, 2).Value = 5230
$trials.Cells.Item(152, 3).Value = 4.688547134399414

$trials.Cells.Item(153, 2).Value = 6228
$trials.Cells.Item(153, 3).Value = 4.747324466705322

$trials.Cells.Item(154, 2).Value = 3491
$trials.Cells.Item(154, 3).Value = 2.116340160369873

$trials.Cells.Item(155, 2).Value = 4557
$trials.Cells.Item(155, 3).Value = 3.36346173286438
$trials.Cells.Item(155, 4).Value = 1

$trials.Cells.Item(156, 2).Value = 1393
$trials.Cells.Item(156, 3).Value = 0.6486704349517822

$trials.Cells.Item(157, 2).Value = 4152
$trials.Cells.Item(157, 3).Value = 2.41569447517395

$trials.Cells.Item(158, 2).Value = 2525
$trials.Cells.Item(158, 3).Value = 0.9495813846588135

$trials.Cells.Item(159, 2).Value = 3591
$trials.Cells.Item(159, 3).Value = 1.869731903076172

$trials.Cells.Item(160, 2).Value = 6862
$trials.Cells.Item(160, 3).Value = 5.98238468170166

$trials.Cells.Item(161, 2).Value = 12851
$trials.Cells.Item(161, 3).Value = 9.533844709396362

$trials.Cells.Item(162, 2).Value = 6858
$trials.Cells.Item(162, 3).Value = 5.150264263153076

$trials.Cells.Item(163, 2).Value = 1879
$trials.Cells.Item(163, 3).Value = 0.9529023170471191

$trials.Cells.Item(164, 2).Value = 3588
$trials.Cells.Item(164, 3).Value = 2.224061489105225

$trials.Cells.Item(165, 2).Value = 7698
$trials.Cells.Item(165, 3).Value = 7.017829895019531

$trials.Cells.Item(166, 2).Value = 2714
$trials.Cells.Item(166, 3).Value = 1.572695255279541

$trials.Cells.Item(167, 2).Value = 320
$trials.Cells.Item(167, 3).Value = 0.1202051639556885

$trials.Cells.Item(168, 2).Value = 6024
$trials.Cells.Item(168, 3).Value = 4.849347114562988

$trials.Cells.Item(169, 2).Value = 4733
$trials.Cells.Item(169, 3).Value = 2.078632831573486

$trials.Cells.Item(170, 2).Value = 6217
$trials.Cells.Item(170, 3).Value = 4.938544273376465

$trials.Cells.Item(171, 2).Value = 7493
$trials.Cells.Item(171, 3).Value = 6.979292869567871

$trials.Cells.Item(172, 2).Value = 6487
$trials.Cells.Item(172, 3).Value = 5.224421739578247

$trials.Cells.Item(173, 2).Value = 2473
$trials.Cells.Item(173, 3).Value = 1.505780935287476

$trials.Cells.Item(174, 2).Value = 1724
$trials.Cells.Item(174, 3).Value = 0.7638120651245117

$trials.Cells.Item(175, 2).Value = 2998
$trials.Cells.Item(175, 3).Value = 1.512568473815918

$trials.Cells.Item(176, 2).Value = 9452
$trials.Cells.Item(176, 3).Value = 9.790473461151123

$trials.Cells.Item(177, 2).Value = 9651
$trials.Cells.Item(177, 3).Value = 8.911376714706421

$trials.Cells.Item(178, 2).Value = 1109
$trials.Cells.Item(178, 3).Value = 0.4606068134307861

$trials.Cells.Item(179, 2).Value = 7344
$trials.Cells.Item(179, 3).Value = 3.806559562683105

$trials.Cells.Item(180, 2).Value = 2797
$trials.Cells.Item(180, 3).Value = 1.332084178924561

$trials.Cells.Item(181, 2).Value = 10338
$trials.Cells.Item(181, 3).Value = 11.88009285926819

$trials.Cells.Item(182, 2).Value = 6353
$trials.Cells.Item(182, 3).Value = 5.309545278549194
$trials.Cells.Item(182, 4).Value = 1

$trials.Cells.Item(183, 2).Value = 8595
$trials.Cells.Item(183, 3).Value = 7.028247356414795

$trials.Cells.Item(184, 2).Value = 3112
$trials.Cells.Item(184, 3).Value = 1.497905015945435

$trials.Cells.Item(185, 2).Value = 3701
$trials.Cells.Item(185, 3).Value = 2.027366161346436

$trials.Cells.Item(186, 2).Value = 4531
$trials.Cells.Item(186, 3).Value = 2.822676658630371

$trials.Cells.Item(187, 2).Value = 550
$trials.Cells.Item(187, 3).Value = 0.2080566883087158

$trials.Cells.Item(188, 2).Value = 9756
$trials.Cells.Item(188, 3).Value = 10.74779772758484

$trials.Cells.Item(189, 2).Value = 8624
$trials.Cells.Item(189, 3).Value = 9.024158239364624

$trials.Cells.Item(190, 2).Value = 7767
$trials.Cells.Item(190, 3).Value = 6.701136827468872

$trials.Cells.Item(191, 2).Value = 13761
$trials.Cells.Item(191, 3).Value = 18.51025700569153

$trials.Cells.Item(192, 2).Value = 2973
$trials.Cells.Item(192, 3).Value = 1.298262596130371

$trials.Cells.Item(193, 2).Value = 3262
$trials.Cells.Item(193, 3).Value = 1.730854034423828

$trials.Cells.Item(194, 2).Value = 2561
$trials.Cells.Item(194, 3).Value = 1.106073617935181

$trials.Cells.Item(195, 2).Value = 3355
$trials.Cells.Item(195, 3).Value = 1.759443759918213

$trials.Cells.Item(196, 2).Value = 6480
$trials.Cells.Item(196, 3).Value = 5.063333034515381

$trials.Cells.Item(197, 2).Value = 6394
$trials.Cells.Item(197, 3).Value = 4.812201023101807
$trials.Cells.Item(197, 4).Value = 1

$trials.Cells.Item(198, 2).Value = 2318
$trials.Cells.Item(198, 3).Value = 0.901400089263916

$trials.Cells.Item(199, 2).Value = 3103
$trials.Cells.Item(199, 3).Value = 1.606194019317627

$trials.Cells.Item(200, 2).Value = 4717
$trials.Cells.Item(200, 3).Value = 2.803488492965698

$trials.Cells.Item(201, 2).Value = 3712
$trials.Cells.Item(201, 3).Value = 1.956302165985107

# Summary sheet row 2
$summary.Cells.Item(2, 1).Value = 160
$summary.Cells.Item(2, 2).Value = 0.06010258197784424
$summary.Cells.Item(2, 3).Value = 0.48

Write-Output "applied"